{"js": "// Apply the Blirix Workshop copy refresh.\n// Each entry is an exact, whole-paragraph-run text match (matchCase, no\n// wildcards) paired with its replacement. We search the whole body for\n// each old string and replace every match - the title string legitimately\n// appears twice (the H1 heading and the bold paragraph near the end) and\n// both copies need the same new text.\nconst replacements = [\n  [\n    \"Play Blirix Workshop Free - Impressive Steampunk Setting\",\n    \"Play Blirix Workshop Free: Impressive Steampunk Adventure\",\n  ],\n  [\n    \"Impressive Steampunk setting\",\n    \"Impressive Steampunk setting creates a fantastic world of adventure\",\n  ],\n  [\n    \"High volatility, infrequent but significant wins\",\n    \"Incredible graphics and animations bring the game to life\",\n  ],\n  [\n    \"Two wild symbols with advanced spells\",\n    \"High volatility offers the chance for significant wins\",\n  ],\n  [\n    \"Free spins with added rewards and advanced spells\",\n    \"Exciting free spins and bonus features add to the gameplay\",\n  ],\n  [\n    \"Limited availability of auto spins\",\n    \"Limited number of paylines\",\n  ],\n  [\n    \"Minimum bet of 0.20 \u20ac may not be suitable for all players\",\n    \"Minimum bet may be higher for some players\",\n  ],\n  [\n    \"Read our Blirix Workshop review and play for free! Enjoy the impressive Steampunk setting, high volatility, and advanced spells during free spins.\",\n    \"Read our review of Blirix Workshop and play for free. Embark on a captivating Steampunk adventure.\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const r of found.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the Blirix Workshop copy refresh via Find & Replace.\n# Each pair is an exact, whole-paragraph-run text match; MatchCase is\n# forced on so e.g. \"Impressive Steampunk Setting\" (title, capital S)\n# and \"Impressive Steampunk setting\" (bullet, lowercase s) don't collide.\n# wdReplaceAll (2) + Wrap:=wdFindContinue (1) sweeps the whole story, so\n# the title string - which legitimately appears twice (H1 heading and the\n# bold paragraph near the end) - gets updated in both places in one call.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"Play Blirix Workshop Free - Impressive Steampunk Setting\", \"Play Blirix Workshop Free: Impressive Steampunk Adventure\"),\n    @(\"Impressive Steampunk setting\", \"Impressive Steampunk setting creates a fantastic world of adventure\"),\n    @(\"High volatility, infrequent but significant wins\", \"Incredible graphics and animations bring the game to life\"),\n    @(\"Two wild symbols with advanced spells\", \"High volatility offers the chance for significant wins\"),\n    @(\"Free spins with added rewards and advanced spells\", \"Exciting free spins and bonus features add to the gameplay\"),\n    @(\"Limited availability of auto spins\", \"Limited number of paylines\"),\n    @(\"Minimum bet of 0.20 \u20ac may not be suitable for all players\", \"Minimum bet may be higher for some players\"),\n    @(\"Read our Blirix Workshop review and play for free! Enjoy the impressive Steampunk setting, high volatility, and advanced spells during free spins.\", \"Read our review of Blirix Workshop and play for free. Embark on a captivating Steampunk adventure.\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
